# Update the "2024" worksheet: a new log entry was added at the top of the
# "September" electricity-log stream (columns R/S) which pushes the existing
# "edc"/"amazeloan"/"hdfc" style log entries down by one row within their
# respective Details/Date column pair, and also pushes the "Group" labels in
# column A down by one row, adding a new "Broadband" entry at row 37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# --- Row 17: new "September" entry (R/S) replaces old "August" entry (P/Q) ---
$ws.Range("P17").Value = ""
$ws.Range("Q17").Value = ""
$ws.Range("R17").Value = "electricity avoid disconnection tangedco"
$ws.Range("S17").Value = "2024-09-03 07:56:18"

# --- Rows 18-21: "edc" August entries shift their date down by one slot ---
$ws.Range("Q18").Value = "2024-08-21 20:15:45"
$ws.Range("Q19").Value = "2024-08-21 20:15:34"
$ws.Range("Q20").Value = "2024-08-21 20:14:29"
$ws.Range("Q21").Value = "2024-08-21 20:12:49"

# --- Row 22: entry moves from July columns (N/O) to August columns (P/Q) ---
$ws.Range("N22").Value = ""
$ws.Range("O22").Value = ""
$ws.Range("P22").Value = "edc"
$ws.Range("Q22").Value = "2024-08-21 20:12:02"

# --- Row 23: "edc" July entry date shifts down by one slot ---
$ws.Range("O23").Value = "2024-07-20 18:52:06"

# --- Row 24: Group label cleared, new "edc" July entry appears ---
$ws.Range("A24").Value = ""
$ws.Range("N24").Value = "edc"
$ws.Range("O24").Value = "2024-07-20 18:51:36"

# --- Rows 25-27: Group labels shift down by one ---
$ws.Range("A25").Value = "OTT"
$ws.Range("A26").Value = "Swiggy"
$ws.Range("A27").Value = "Others"

# --- Row 27: old "broker" entry (R/S) removed (moved to row 28) ---
$ws.Range("R27").Value = ""
$ws.Range("S27").Value = ""

# --- Row 28: "amazeloan" entry becomes the "broker" entry from row 27 ---
$ws.Range("R28").Value = "broker"
$ws.Range("S28").Value = "2024-09-01 22:35:38"

# --- Rows 29-31: "amazeloan" September entries shift date down by one slot ---
$ws.Range("S29").Value = "2024-09-01 10:12:03"
$ws.Range("S30").Value = "2024-09-01 09:42:38"
$ws.Range("S31").Value = "2024-09-01 09:29:24"

# --- Row 32: entry moves from August columns (P/Q) to September columns (R/S) ---
$ws.Range("P32").Value = ""
$ws.Range("Q32").Value = ""
$ws.Range("R32").Value = "amazeloan"
$ws.Range("S32").Value = "2024-09-01 09:27:06"

# --- Rows 33-35: "hdfc" August entries shift date down by one slot ---
$ws.Range("Q33").Value = "2024-08-30 12:15:48"
$ws.Range("Q34").Value = "2024-08-21 20:17:10"
$ws.Range("Q35").Value = "2024-08-21 20:16:45"

# --- Row 36: Group label cleared, new "hdfc" August entry appears ---
$ws.Range("A36").Value = ""
$ws.Range("P36").Value = "hdfc"
$ws.Range("Q36").Value = "2024-08-21 20:15:50"

# --- Row 37 (new row): Group label "Broadband" ---
$ws.Range("A37").Value = "Broadband"
